# fix #140 utterance true turn to boolean
#
# Adds a new "AMAZON.YesIntent" intent:
#  - INTENT sheet: new row with the intent name in column A
#  - UTTERANCES_MAIN sheet: new column with the intent name as header and
#    its sample utterances ("true" / "yes") below it. The "true" utterance
#    is explicitly formatted as Text so Excel/LibreOffice doesn't coerce it
#    into the boolean TRUE value (that's the bug being fixed).

$wb = $excel.ActiveWorkbook

$intentSheet = $wb.Worksheets.Item("INTENT")
$utterancesSheet = $wb.Worksheets.Item("UTTERANCES_MAIN")

# --- INTENT sheet: append the new intent name in column A ---
$introw = $intentSheet.Cells.Item(10, 1)
$introw.Value = "AMAZON.YesIntent"

# --- UTTERANCES_MAIN sheet: add a new column (F) for this intent ---
$utterancesSheet.Cells.Item(1, 6).Value = "AMAZON.YesIntent"

# The literal word "true" must stay text, not become a real boolean value.
$trueCell = $utterancesSheet.Cells.Item(2, 6)
$trueCell.NumberFormat = "@"
$trueCell.Value = "true"

$utterancesSheet.Cells.Item(3, 6).Value = "yes"
